$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 3923
$ws.Range("B2").Value = 5069
$ws.Range("C2").Value = 8992
$ws.Range("D2").Value = 3159.04890750022
$ws.Range("E2").Value = -1493.329938260623
$ws.Range("F2").Value = 0.3157816541360874
$ws.Range("G2").Value = -0.1484917256038622
$ws.Range("H2").Value = 0.05405992794593976
$ws.Range("I2").Value = 12392948.86412337
$ws.Range("J2").Value = -7569689.457043095
$ws.Range("K2").Value = 0.4362766903914591
$ws.Range("L2").Value = 2.115439345694607
$ws.Range("M2").Value = 1.637180618102181
$ws.Range("N2").Value = 4823259.40708028
